# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Price (D) and Volume(1h) (E) columns are stored as text in the sheet, so
# NumberFormat is forced to "@" before each Price write to stop Excel's
# automatic type-inference from silently converting plain-decimal strings
# (e.g. "0.4678") into numeric cells and mangling their exact text form.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.228.41'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.861.20'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.53'
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4678'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2870'
$ws.Range("E8").Value = '  +0.93%  '
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.69'
$ws.Range("E10").Value = '  +3.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07932'
$ws.Range("E11").Value = '  +0.32%  '
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.863.01'
$ws.Range("E13").Value = '  -0.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.174'
$ws.Range("E14").Value = '  +0.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6809'
$ws.Range("E15").Value = '  +0.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '267.78'
$ws.Range("E16").Value = '  -6.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.232.40'
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.71'
$ws.Range("E18").Value = '  +7.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007404'
$ws.Range("E20").Value = '  +1.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.111.14'
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.327'
$ws.Range("E22").Value = '  -4.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.191'
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.33'
$ws.Range("E25").Value = '  +0.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.229'
$ws.Range("E26").Value = '  -1.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.91'
$ws.Range("E27").Value = '  -1.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.966'
$ws.Range("E28").Value = '  +2.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.385'
$ws.Range("E29").Value = '  +1.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09841'
$ws.Range("E30").Value = '  +1.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.394'
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.474'
$ws.Range("E32").Value = '  +0.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.074'
$ws.Range("E33").Value = '  -0.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04699'
$ws.Range("E34").Value = '  -0.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.134'
$ws.Range("E35").Value = '  +0.94%  '
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("E37").Value = '  -0.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01878'
$ws.Range("E38").Value = '  +0.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.638'
$ws.Range("E39").Value = '  +3.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.250'
$ws.Range("E40").Value = '  -2.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.44'
$ws.Range("E41").Value = '  +0.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.939'
$ws.Range("E42").Value = '  -0.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8460'
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4173'
$ws.Range("E44").Value = '  -0.39%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.40'
$ws.Range("E46").Value = '  -0.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '959.45'
$ws.Range("E47").Value = '  +3.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.164'
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.242'
$ws.Range("E49").Value = '  -0.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.14'
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05657'
$ws.Range("E51").Value = '  +0.42%  '
